$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.224.29"

$ws.Range("D3").Value = "1.789.64"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "227.05"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "32.26"
$ws.Range("E8").Value = "  -1.18%  "

$ws.Range("E9").Value = "  +2.72%  "

$ws.Range("E10").Value = "  -2.66%  "

$ws.Range("E11").Value = "  +0.85%  "

$ws.Range("D12").Value = "2.047.56"
$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").Value = "11.35"
$ws.Range("E13").Value = "  +3.97%  "

$ws.Range("D14").Value = "1.797.25"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").Value = "34.134.50"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("E17").Value = "  +2.00%  "

$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "245.70"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("D20").Value = "0.0₃0781"
$ws.Range("E20").Value = "  -0.41%  "

$ws.Range("D21").Value = "10.93"
$ws.Range("E21").Value = "  +2.13%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "4.12"
$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("E24").Value = "  -1.14%  "

$ws.Range("D25").Value = "161.70"
$ws.Range("E25").Value = "  +0.98%  "

$ws.Range("E26").Value = "  +1.85%  "

$ws.Range("D27").Value = "16.31"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("E30").Value = "  +1.03%  "

$ws.Range("E31").Value = "  +1.81%  "

$ws.Range("E32").Value = "  +1.14%  "

$ws.Range("E33").Value = "  +3.30%  "

$ws.Range("E34").Value = "  +1.39%  "

$ws.Range("D35").Value = "1.442.76"
$ws.Range("E35").Value = "  +3.94%  "

$ws.Range("D36").Value = "0.648"
$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("D37").Value = "2.43"
$ws.Range("E37").Value = "  +10.68%  "

$ws.Range("E38").Value = "  +2.87%  "

$ws.Range("E39").Value = "  -0.70%  "

$ws.Range("D40").Value = "80.54"
$ws.Range("E40").Value = "  +3.57%  "

$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("E42").Value = "  +1.08%  "

$ws.Range("E43").Value = "  +0.56%  "

$ws.Range("D44").Value = "13.32"
$ws.Range("E44").Value = "  +5.09%  "

$ws.Range("E45").Value = "  +2.58%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0139"
$ws.Range("E46").Value = "  +0.13%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "6.07"
$ws.Range("E47").Value = "  +4.59%  "

$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("D49").Value = "107.78"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").Value = "1.949.07"
$ws.Range("E50").Value = "  +0.76%  "

$ws.Range("E51").Value = "  +0.08%  "
